$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the start date for the "Pre-emergency" row (A2) from 2020-03-01 to 2020-01-01
$ws.Range("A2").Value = "2020-01-01"

# Update the active selection to F11
$ws.Range("F11").Select()
